# Generate Report for Handback
# Re-run of the handback status report: the two tracked e2e markdown files
# ("b9e5070f-...md" and "da3eb767-...md") were regenerated under new GUID
# names ("aa50d49d-...md" and "ffff44e0b200-...md") with refreshed
# handoff/handback timestamps and xliff artifact names.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "b9e5070f-5640-448c-82e6-5112491c6b28"
$newGuid1 = "aa50d49d-f69f-41b0-9431-6339ccc40d93"
$oldGuid2 = "da3eb767-701e-44f3-b29b-7d92dd12abb6"
$newGuid2 = "ffff44e0b200-f768-4d38-b224-7096d2caa4ad"

$oldHash1 = "285f8a83521a15c33e4c2d80c83ae0b1ab42b93f"
$newHash1 = "8e2a1fd2ea3c3e0d352ac44c978a8cc0e4dc7e2d"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("G2").Value = "2016-08-21 09:06:50"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("G3").Value = "2016-08-21 09:06:50"

# Keep the hyperlink display text in sync with the new cell text.
$wsOverview.Range("B2").Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid1.md"
$wsOverview.Range("B3").Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid2.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid1.md"
$wsZhCn.Range("G2").Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-21 09:06:46"
$wsZhCn.Range("I2").Value = "$newGuid1.md"
$wsZhCn.Range("J2").Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-21 09:07:06"

$wsZhCn.Range("A3").Value = "$newGuid2.md"
$wsZhCn.Range("G3").Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-21 09:06:46"
$wsZhCn.Range("I3").Value = "$newGuid2.md"
$wsZhCn.Range("J3").Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-21 09:07:06"

$wsZhCn.Range("A2").Hyperlinks.Item(1).TextToDisplay = "$newGuid1.md"
$wsZhCn.Range("I2").Hyperlinks.Item(1).TextToDisplay = "$newGuid1.md"
$wsZhCn.Range("A3").Hyperlinks.Item(1).TextToDisplay = "$newGuid2.md"
$wsZhCn.Range("I3").Hyperlinks.Item(1).TextToDisplay = "$newGuid2.md"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid1.md"
$wsDeDe.Range("G2").Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-21 09:06:50"
$wsDeDe.Range("I2").Value = "$newGuid1.md"
$wsDeDe.Range("J2").Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-21 09:07:12"

$wsDeDe.Range("A3").Value = "$newGuid2.md"
$wsDeDe.Range("G3").Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-21 09:06:50"
$wsDeDe.Range("I3").Value = "$newGuid2.md"
$wsDeDe.Range("J3").Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-21 09:07:12"

$wsDeDe.Range("A2").Hyperlinks.Item(1).TextToDisplay = "$newGuid1.md"
$wsDeDe.Range("I2").Hyperlinks.Item(1).TextToDisplay = "$newGuid1.md"
$wsDeDe.Range("A3").Hyperlinks.Item(1).TextToDisplay = "$newGuid2.md"
$wsDeDe.Range("I3").Hyperlinks.Item(1).TextToDisplay = "$newGuid2.md"
